# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the existing row 75 ("SBS Main Indicators...")
# so everything from the old row 75 downward shifts down by one. This creates
# the new empty row 75 and pushes the remaining source/footnote rows into
# their new positions (76-78, 81-86) while preserving their existing styles.
$ws.Rows("75").Insert()

# Remove the hyperlink that used to sit on (old) A76 -> (now) A77, and plain-text it
# (restyle it to match the other plain "source" footnote rows: italic, no underline/color).
$ws.Range("A77").Hyperlinks.Delete()
$ws.Range("A77").Value = ""
$ws.Range("A77").Font.Italic = $true
$ws.Range("A77").Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleNone
$ws.Range("A77").Font.ColorIndex = [Microsoft.Office.Interop.Excel.XlColorIndex]::xlColorIndexAutomatic

# The eurostat URL moves from A77 down into A78, now as plain "source" text.
$ws.Range("A78").Value = "http://epp.eurostat.ec.europa.eu/portal/page/portal/european_business/data/database"

# Old row 81 (the removed "Provisional results..." footnote) becomes a second
# "Swiss Statistics" label row, styled like the other "Swiss Statistics" rows.
$ws.Range("A82").Value = "Swiss Statistics"

# Old row 83 ("Business indicators...") is replaced with the new citation text
# about the SME definition.
$ws.Range("A84").Value = [char]0x201C + "Small and medium-sized enterprises" + [char]0x201D + "  in Swiss Federal Statistical Office > Data Library > Definitions."

# Old row 85 (the removed "Sructural Business Statistics..." footnote) becomes
# a second "SBS Eurostat" label row.
$ws.Range("A86").Value = "SBS Eurostat"
